# Update the topic speaker on the "az900" sheet.
# Row 2 (the "describe-core-architectural-components-of-azure" topic) had
# "Sowmya" listed as the speaker/name in column A; it is now "saikiran".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("az900")
$ws.Range("A2").Value = "saikiran"
